# Weekly data refresh: insert the newest week's "Limón" price rows at the
# top of the data block (row 844), pushing the rest of the table down by
# two rows. This mirrors how the upstream feed prepends the latest report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 844-845; everything currently at 844.. shifts
# down to 846.. (Excel copies the formatting of the row above into the
# freshly inserted rows, same as a manual Insert in the UI).
$ws.Rows("844:845").Insert()

# Row 844 — "1a amarillo" grade, newest report date.
$ws.Range("A844").Value2 = 11
$ws.Range("B844").Value2 = "Vega Monumental Concepción"
$ws.Range("C844").Value2 = "Bíobío"
$ws.Range("D844").Value2 = 45223
$ws.Range("E844").Value2 = 8
$ws.Range("F844").Value2 = "Fruta"
$ws.Range("G844").Value2 = 100102
$ws.Range("H844").Value2 = "Cítricos"
$ws.Range("I844").Value2 = 100102003
$ws.Range("J844").Value2 = "Limón"
$ws.Range("K844").Value2 = "Sin especificar"
$ws.Range("L844").Value2 = "1a amarillo"
$ws.Range("M844").Value2 = 500
$ws.Range("N844").Value2 = 9000
$ws.Range("O844").Value2 = 9500
$ws.Range("P844").Value2 = 9300
$ws.Range("Q844").Value2 = "$/malla 18 kilos"
$ws.Range("R844").Value2 = "Provincia de Melipilla"
$ws.Range("S844").Value2 = 517
$ws.Range("T844").Value2 = 18

# Row 845 — "2a amarillo" grade, same newest report date.
$ws.Range("A845").Value2 = 11
$ws.Range("B845").Value2 = "Vega Monumental Concepción"
$ws.Range("C845").Value2 = "Bíobío"
$ws.Range("D845").Value2 = 45223
$ws.Range("E845").Value2 = 8
$ws.Range("F845").Value2 = "Fruta"
$ws.Range("G845").Value2 = 100102
$ws.Range("H845").Value2 = "Cítricos"
$ws.Range("I845").Value2 = 100102003
$ws.Range("J845").Value2 = "Limón"
$ws.Range("K845").Value2 = "Sin especificar"
$ws.Range("L845").Value2 = "2a amarillo"
$ws.Range("M845").Value2 = 200
$ws.Range("N845").Value2 = 7500
$ws.Range("O845").Value2 = 7500
$ws.Range("P845").Value2 = 7500
$ws.Range("Q845").Value2 = "$/malla 18 kilos"
$ws.Range("R845").Value2 = "Provincia de Melipilla"
$ws.Range("S845").Value2 = 417
$ws.Range("T845").Value2 = 18
